$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Coin/Link columns remain plain text so that values
# like '1.00' or '0.0000133' are not coerced into numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.215.03'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.601.56'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.02'
$ws.Range('E5').Value = '  +4.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.04'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.45'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('E10').Value = '  +1.96%  '
$ws.Range('E11').Value = '  +1.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.067.17'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '59.157.18'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.53'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.635.55'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '343.21'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.11'
$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.52'
$ws.Range('E23').Value = '  +2.32%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Kaspa'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.165'
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.407'
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.19'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0736'
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.70'
$ws.Range('E30').Value = '  +8.93%  '
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.70'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '149.58'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '37.15'
$ws.Range('E35').Value = '  +2.41%  '
$ws.Range('E36').Value = '  -1.34%  '
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.833'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.812'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.55'
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '274.98'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0955'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.940.86'
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.41'
$ws.Range('E49').Value = '  +2.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.50'
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.10'
$ws.Range('E51').Value = '  -1.78%  '
